$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" (column E) values for rows 16-22 and their corresponding
# "Valor Mora" date-like numeric values (column F) get reversed top-to-bottom,
# i.e. the previous account-statement periods are removed and new ones (most
# recent first) are entered in their place.
$ws.Range("E16").Value = "2311"
$ws.Range("E17").Value = "2310"
$ws.Range("E18").Value = "2309"
$ws.Range("E19").Value = "2308"
$ws.Range("E20").Value = "2307"
$ws.Range("E21").Value = "2306"
$ws.Range("E22").Value = "2305"

$ws.Range("F16").Value = 33918
$ws.Range("F17").Value = 36341
$ws.Range("F18").Value = 36341
$ws.Range("F19").Value = 36341
$ws.Range("F20").Value = 36341
$ws.Range("F21").Value = 36341
$ws.Range("F22").Value = 36341
